{"js": "// Replace the division-problem text in each table cell with its new value.\n// Each old value is unique in the document, so a targeted search+replace\n// per pair reliably retargets exactly the intended run without disturbing\n// any other content (e.g. the date heading paragraph).\nconst replacements = [\n  [\"841\u00f74=\", \"236\u00f75=\"],\n  [\"289\u00f77=\", \"537\u00f74=\"],\n  [\"623\u00f76=\", \"413\u00f77=\"],\n  [\"743\u00f79=\", \"238\u00f78=\"],\n  [\"873\u00f78=\", \"968\u00f72=\"],\n  [\"500\u00f77=\", \"621\u00f78=\"],\n  [\"231\u00f78=\", \"154\u00f74=\"],\n  [\"402\u00f78=\", \"376\u00f78=\"],\n  [\"364\u00f73=\", \"401\u00f77=\"],\n  [\"801\u00f76=\", \"615\u00f73=\"],\n  [\"458\u00f74=\", \"347\u00f72=\"],\n  [\"740\u00f74=\", \"172\u00f73=\"],\n  [\"803\u00f76=\", \"343\u00f79=\"],\n  [\"471\u00f74=\", \"695\u00f79=\"],\n  [\"425\u00f72=\", \"156\u00f75=\"],\n  [\"715\u00f79=\", \"888\u00f76=\"],\n  [\"424\u00f77=\", \"182\u00f78=\"],\n  [\"704\u00f78=\", \"673\u00f77=\"],\n  [\"686\u00f74=\", \"106\u00f78=\"],\n  [\"162\u00f78=\", \"975\u00f79=\"],\n  [\"133\u00f72=\", \"406\u00f73=\"],\n  [\"821\u00f73=\", \"420\u00f76=\"],\n  [\"963\u00f78=\", \"408\u00f76=\"],\n  [\"229\u00f73=\", \"114\u00f72=\"],\n  [\"196\u00f78=\", \"193\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in each table cell with its new value.\n# Each old value is unique across the document, so Find/Replace against the\n# whole document's content range safely retargets exactly the intended run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"841\u00f74=\", \"236\u00f75=\"),\n    @(\"289\u00f77=\", \"537\u00f74=\"),\n    @(\"623\u00f76=\", \"413\u00f77=\"),\n    @(\"743\u00f79=\", \"238\u00f78=\"),\n    @(\"873\u00f78=\", \"968\u00f72=\"),\n    @(\"500\u00f77=\", \"621\u00f78=\"),\n    @(\"231\u00f78=\", \"154\u00f74=\"),\n    @(\"402\u00f78=\", \"376\u00f78=\"),\n    @(\"364\u00f73=\", \"401\u00f77=\"),\n    @(\"801\u00f76=\", \"615\u00f73=\"),\n    @(\"458\u00f74=\", \"347\u00f72=\"),\n    @(\"740\u00f74=\", \"172\u00f73=\"),\n    @(\"803\u00f76=\", \"343\u00f79=\"),\n    @(\"471\u00f74=\", \"695\u00f79=\"),\n    @(\"425\u00f72=\", \"156\u00f75=\"),\n    @(\"715\u00f79=\", \"888\u00f76=\"),\n    @(\"424\u00f77=\", \"182\u00f78=\"),\n    @(\"704\u00f78=\", \"673\u00f77=\"),\n    @(\"686\u00f74=\", \"106\u00f78=\"),\n    @(\"162\u00f78=\", \"975\u00f79=\"),\n    @(\"133\u00f72=\", \"406\u00f73=\"),\n    @(\"821\u00f73=\", \"420\u00f76=\"),\n    @(\"963\u00f78=\", \"408\u00f76=\"),\n    @(\"229\u00f73=\", \"114\u00f72=\"),\n    @(\"196\u00f78=\", \"193\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
